$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ticket rows appended after the existing data (rows 1-136).
$data = @(
    @("2024-05-21", "12:03:58", "Etiquetadora",  "-", "-", "-", "-", "12:06:11", "0:02:13"),
    @("2024-05-21", "12:09:37", "Etiquetadora2", "-", "-", "-", "-", "12:09:39", "0:00:02"),
    @("2024-05-21", "12:10:19", "Etiquetadora",  "-", "-", "-", "-", "12:10:20", "0:00:01")
)

$startRow = 137
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]

    # Column A holds a yyyy-mm-dd looking string ("2024-05-21"). A plain
    # .Value assignment would get auto-parsed into a date serial number by
    # Excel, which does not match the source data (plain text). Force the
    # cell to text mode for the assignment, then restore the cell to the
    # workbook's normal (General) style so no stray number format lingers
    # on the cell itself.
    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $values[0]
    $dateCell.Style = "Normal"

    for ($c = 1; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
